$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 227, shifting existing rows 227:241 down to 228:242
$ws.Range("A227").EntireRow.Insert()

# Populate the newly inserted row 227 with the new weekly record
$ws.Cells.Item(227, 1).Value = 4
$ws.Cells.Item(227, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(227, 3).Value = "Los Lagos"
$ws.Cells.Item(227, 4).Value = 44714
$ws.Cells.Item(227, 4).NumberFormat = $ws.Cells.Item(228, 4).NumberFormat
$ws.Cells.Item(227, 5).Value = 10
$ws.Cells.Item(227, 6).Value = 100112032
$ws.Cells.Item(227, 7).Value = "Zapallo italiano"
$ws.Cells.Item(227, 8).Value = "Sin especificar"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 70
$ws.Cells.Item(227, 11).Value = 15000
$ws.Cells.Item(227, 12).Value = 17000
$ws.Cells.Item(227, 13).Value = 16000
$ws.Cells.Item(227, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(227, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(227, 16).Value = 320
$ws.Cells.Item(227, 17).Value = 50
$ws.Cells.Item(227, 18).Value = "Hortaliza"
